# Applies the Gaussian Quadrature Scheme related edits to the Averaged Intensities workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab/title.
$ws.Name = "GossF"

# Correct tiny floating point differences in row 13 (recomputed with updated quadrature scheme).
$ws.Range("C13").Value = 1.005877367761607
$ws.Range("F13").Value = 1.005877367761607
$ws.Range("L13").Value = 0.9978877434310979
$ws.Range("M13").Value = 0.9925594138841238

# Correct tiny floating point differences in row 15.
$ws.Range("C15").Value = 0.9841436777953232
$ws.Range("F15").Value = 0.9841436777953232

# Add a new row 16 of averaged intensity data (HKL index 14, HexGrid-60degTilt5degRes scheme).
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.634508309208395
$ws.Range("D16").Value = 2.056121120826005
$ws.Range("E16").Value = 1.021648733139336
$ws.Range("F16").Value = 1.634508309208395
$ws.Range("G16").Value = 0.7277282058053257
$ws.Range("H16").Value = 2.018785056086185
$ws.Range("I16").Value = 0.7717809265185526
$ws.Range("J16").Value = 2.056121120826005
$ws.Range("K16").Value = 1.53888492698267
$ws.Range("L16").Value = 1.586696618095533
$ws.Range("M16").Value = 1.3717620585973

# Apply the same style as column A in the other data rows (bordered, bold, centered header-like style).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
